# Append a new data row (2025/10/01, 水, 17, 3) to the bottom of the
# existing table on the active sheet, mirroring the style of the prior
# data rows (no explicit cell style; date column stored as plain text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$dateCell = $ws.Cells.Item($newRow, 1)
# Force text storage so the date-like string isn't auto-converted into a
# serial date number, then drop back to the default "Normal" style so the
# new row doesn't pick up a stray number-format style (matching the
# un-styled cells used by the rest of the data rows).
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025/10/01"
$dateCell.Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = "水"
$ws.Cells.Item($newRow, 3).Value = 17
$ws.Cells.Item($newRow, 4).Value = 3
